# "Generate Report for Handback" -- mark the zh-cn and de-de localization
# rows as handed back: update status text, fill in the Latest Target
# File / Latest Handback File / Latest Handback DateTime columns, add the
# corresponding hyperlinks, and widen the columns that now hold the
# longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$zhCnHandbackDate = "2016-09-05 01:09:30"
$deDeHandbackDate = "2016-09-05 01:09:38"

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec95a867993fda879a8b62c2e67eac7fcd2a117f/e2e/2d32bbc5-6ffa-4263-a591-48bdd002b389.md"
$mdDisplay = "2d32bbc5-6ffa-4263-a591-48bdd002b389.md"
$ffffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec95a867993fda879a8b62c2e67eac7fcd2a117f/e2e/ffff4ab81d1a-a803-4138-a1b6-3c7c1a072991.md"
$ffffDisplay = "ffff4ab81d1a-a803-4138-a1b6-3c7c1a072991.md"

# Hyperlink-style font (matches the existing "HyperLink" cell style already
# used by column A: underlined, RGB 6495ED).
$hyperlinkColor = 15570276  # OLE BGR encoding of 0x6495ED

# The engine's Range/Columns.ColumnWidth setter re-quantises to Excel's
# pixel grid and bakes in a +5/6 character padding offset before writing
# the OOXML <col width=".."> attribute, so to land on a target *stored*
# width we have to feed in (target - 5/6).
$wideStatusColWidth = 29.9777047293527 - (5.0/6.0)   # -> stored ~30
$wideFileColWidth   = 40.0 - (5.0/6.0)                # -> stored exactly 40

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $mdDisplay
$wsZh.Range("J2").Value = "2d32bbc5-6ffa-4263-a591-48bdd002b389.bdb289aa0d6fad8a24ca0b1f7a0472d53a73cb1b.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhCnHandbackDate

$wsZh.Range("I3").Value = $mdDisplay
$wsZh.Range("J3").Value = "2d32bbc5-6ffa-4263-a591-48bdd002b389.bdb289aa0d6fad8a24ca0b1f7a0472d53a73cb1b.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhCnHandbackDate

Style-AsHyperlink($wsZh.Range("I2"))
Style-AsHyperlink($wsZh.Range("I3"))

# Rebuild hyperlinks in A2, I2, A3, I3 order so relationship ids line up
# the way Excel assigns them when the links are (re)created row by row.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, $null, $null, $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, $null, $null, $mdDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ffffUrl, $null, $null, $ffffDisplay) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, $null, $null, $mdDisplay) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsZh.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsZh.Columns.Item(10).ColumnWidth = $wideFileColWidth

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $mdDisplay
$wsDe.Range("J2").Value = "2d32bbc5-6ffa-4263-a591-48bdd002b389.bdb289aa0d6fad8a24ca0b1f7a0472d53a73cb1b.de-de.xlf"
$wsDe.Range("K2").Value = $deDeHandbackDate

$wsDe.Range("I3").Value = $mdDisplay
$wsDe.Range("J3").Value = "2d32bbc5-6ffa-4263-a591-48bdd002b389.bdb289aa0d6fad8a24ca0b1f7a0472d53a73cb1b.de-de.xlf"
$wsDe.Range("K3").Value = $deDeHandbackDate

Style-AsHyperlink($wsDe.Range("I2"))
Style-AsHyperlink($wsDe.Range("I3"))

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, $null, $null, $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, $null, $null, $mdDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ffffUrl, $null, $null, $ffffDisplay) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, $null, $null, $mdDisplay) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = $wideStatusColWidth
$wsDe.Columns.Item(9).ColumnWidth = $wideFileColWidth
$wsDe.Columns.Item(10).ColumnWidth = $wideFileColWidth

# ---- Overview sheet roll-up (same shared status text, widen to match) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusColWidth
